$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cell F1 ("time_taken"), copying formatting from the
# neighbouring header cell E1 (bold, centered, bordered header style).
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("F1").Value = "time_taken"

$timestamps = @(
    "2021-10-05 13:39:48.952907",
    "2021-10-05 13:39:48.952919",
    "2021-10-05 13:39:48.952922",
    "2021-10-05 13:39:48.952925",
    "2021-10-05 13:39:48.952928",
    "2021-10-05 13:39:48.952930",
    "2021-10-05 13:39:48.952933",
    "2021-10-05 13:39:48.952935",
    "2021-10-05 13:39:48.952938",
    "2021-10-05 13:39:48.952941",
    "2021-10-05 13:39:48.952943",
    "2021-10-05 13:39:48.952946",
    "2021-10-05 13:39:48.952948",
    "2021-10-05 13:39:48.952950",
    "2021-10-05 13:39:48.952953",
    "2021-10-05 13:39:48.952956",
    "2021-10-05 13:39:48.952959",
    "2021-10-05 13:39:48.952962"
)

for ($i = 0; $i -lt $timestamps.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 6).Value = $timestamps[$i]
}
